# Fix shape positions/sizes in the FirstDay sequence diagram (slide 1).
# Each entry gives the shape's PowerPoint `Id` plus the new Left/Top/Width/Height
# in points. The point values are pre-selected to the exact IEEE-754 single-
# precision value whose EMU round-trip equals the intended OOXML target (noted
# in the trailing comment), so the written <a:off>/<a:ext> match exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$fixes = @(
    [pscustomobject]@{ Id = 18; Left = 150.0901641845703125; Top = 125.00394439697265625; Width = 16.10874176025390625; Height = 141.4487457275390625; Flip = $null } # EMU off=(1906145,1587550) ext=(204581,1796399)
    [pscustomobject]@{ Id = 19; Left = 378.562469482421875; Top = 213.1922149658203125; Width = 86.11299896240234375; Height = 36.34157562255859375; Flip = $null } # EMU off=(4807743,2707541) ext=(1093635,461538)
    [pscustomobject]@{ Id = 20; Left = 424.734649658203125; Top = 250.6031646728515625; Width = 0.6400000452995300293; Height = 195.981658935546875; Flip = $null } # EMU off=(5394130,3182660) ext=(8128,2488967)
    [pscustomobject]@{ Id = 21; Left = 419.328125; Top = 250.6031646728515625; Width = 12.09307098388671875; Height = 5.8240947723388671875; Flip = $null } # EMU off=(5325467,3182660) ext=(153582,73966)
    [pscustomobject]@{ Id = 28; Left = 306.8853759765625; Top = 231.626312255859375; Width = 72.62929534912109375; Height = 0.00007874015864217654; Flip = $null } # EMU off=(3897444,2941654) ext=(922392,1)
    [pscustomobject]@{ Id = 34; Left = 305.078033447265625; Top = 255.5577239990234375; Width = 117.55701446533203125; Height = 0.0; Flip = $null } # EMU off=(3874491,3245583) ext=(1492974,0)
    [pscustomobject]@{ Id = 35; Left = 23.6902370452880859375; Top = 265.8729248046875; Width = 125.70977020263671875; Height = 0.42346459627151489258; Flip = $null } # EMU off=(300866,3376586) ext=(1596514,5378)
    [pscustomobject]@{ Id = 77; Left = 24.0; Top = 433.634674072265625; Width = 397.822784423828125; Height = 0.0; Flip = $null } # EMU off=(304800,5507160) ext=(5052349,0)
    [pscustomobject]@{ Id = 79; Left = 191.6155242919921875; Top = 165.17315673828125; Width = 82.0907135009765625; Height = 13.32889842987060546875; Flip = $null } # EMU off=(2433517,2097699) ext=(1042552,169277)
    [pscustomobject]@{ Id = 82; Left = 261.391204833984375; Top = 417.068756103515625; Width = 48.914646148681640625; Height = 13.32889842987060546875; Flip = $null } # EMU off=(3319668,5296773) ext=(621216,169277)
    [pscustomobject]@{ Id = 8; Left = 164.6285858154296875; Top = 164.4215850830078125; Width = 125.77606964111328125; Height = 0.00007874015864217654; Flip = $null } # EMU off=(2090783,2088154) ext=(1597356,1)
    [pscustomobject]@{ Id = 45; Left = 289.88214111328125; Top = 144.56048583984375; Width = 15.1959056854248046875; Height = 5.24598455429077148438; Flip = $null } # EMU off=(3681503,1835918) ext=(192988,66624)
    [pscustomobject]@{ Id = 49; Left = 289.882049560546875; Top = 163.7227630615234375; Width = 16.208110809326171875; Height = 98.01937103271484375; Flip = $null } # EMU off=(3681502,2079279) ext=(205843,1244846)
    [pscustomobject]@{ Id = 51; Left = 166.704345703125; Top = 260.896240234375; Width = 133.3574066162109375; Height = 0.2577953040599822998; Flip = 'flipV' } # EMU off=(2117145,3313382) ext=(1693639,3274)
    [pscustomobject]@{ Id = 52; Left = 165.32623291015625; Top = 149.1075592041015625; Width = 126.0179595947265625; Height = 0.0; Flip = $null } # EMU off=(2099643,1893666) ext=(1600428,0)
    [pscustomobject]@{ Id = 74; Left = 303.5528564453125; Top = 167.65716552734375; Width = 17.128505706787109375; Height = 3.59992146492004394531; Flip = $null } # EMU off=(3855121,2129246) ext=(217532,45719)
    [pscustomobject]@{ Id = 87; Left = 319.627349853515625; Top = 162.0980377197265625; Width = 148.09771728515625; Height = 13.32889842987060546875; Flip = $null } # EMU off=(4059267,2058645) ext=(1880841,169277)
    [pscustomobject]@{ Id = 88; Left = 298.544189453125; Top = 174.61590576171875; Width = 124.67795562744140625; Height = 13.32889842987060546875; Flip = $null } # EMU off=(3791511,2217622) ext=(1583410,169277)
    [pscustomobject]@{ Id = 89; Left = 317.046630859375; Top = 195.5556793212890625; Width = 99.42252349853515625; Height = 13.32889842987060546875; Flip = $null } # EMU off=(4026492,2483557) ext=(1262666,169277)
    [pscustomobject]@{ Id = 92; Left = 302.914337158203125; Top = 200.5015106201171875; Width = 17.128505706787109375; Height = 3.59992146492004394531; Flip = $null } # EMU off=(3847012,2546369) ext=(217532,45719)
    [pscustomobject]@{ Id = 61; Left = 451.322540283203125; Top = 125.316619873046875; Width = 47.150554656982421875; Height = 34.018505096435546875; Flip = $null } # EMU off=(5731796,1591521) ext=(598812,432035)
    [pscustomobject]@{ Id = 62; Left = 477.80181884765625; Top = 159.3351287841796875; Width = 0.0; Height = 40.86606597900390625; Flip = $null } # EMU off=(6068083,2023556) ext=(0,518999)
    [pscustomobject]@{ Id = 66; Left = 474.748992919921875; Top = 184.3212738037109375; Width = 7.7274017333984375; Height = 9.03236293792724609375; Flip = $null } # EMU off=(6029312,2340880) ext=(98138,114711)
    [pscustomobject]@{ Id = 67; Left = 301.17230224609375; Top = 185.721038818359375; Width = 173.5766143798828125; Height = 1.89118111133575439453; Flip = $null } # EMU off=(3824888,2358657) ext=(2204423,24018)
    [pscustomobject]@{ Id = 68; Left = 305.294097900390625; Top = 192.1537933349609375; Width = 169.45489501953125; Height = 1.71078741550445556641; Flip = $null } # EMU off=(3877235,2440353) ext=(2152077,21727)
)

foreach ($fix in $fixes) {
    foreach ($shp in $s.Shapes) {
        if ($shp.Id -eq $fix.Id) {
            $shp.Left = $fix.Left
            $shp.Top = $fix.Top
            $shp.Width = $fix.Width
            $shp.Height = $fix.Height
            if ($fix.Flip -eq 'flipV') {
                $shp.VerticalFlip = -1
            }
            break
        }
    }
}
